# mockportfolio.xlsx update — "Add files via upload"
# Renames the sheet, refreshes the OXY row with real trade data, reworks the
# IRR/"months" footnote area into a Portfolio-vs-SPY retrospective analysis
# block, and appends a couple of closing commentary lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename the sheet
# ------------------------------------------------------------------
$ws.Name = "Main"

# ------------------------------------------------------------------
# 2) OXY row (row 10) — fill in the real buy/current numbers that used to
#    be blank (was producing #DIV/0! in F10/J10, and G10 had no formula).
# ------------------------------------------------------------------
$ws.Range("B10").Value = 196
$ws.Range("C10").Value = 170
$ws.Range("E10").Value = 45784
$ws.Range("I10").Value = 211
$ws.Range("G10").Formula = "=C10/B10-1"

# ------------------------------------------------------------------
# 3) Row 12 header strip & shared-formula range for column G both
#    auto-follow once the string table is renumbered — nothing else to do
#    there. Extend the G column shared % formula down through row 10.
# ------------------------------------------------------------------
$ws.Range("G3:G10").FormulaR1C1 = "=RC3/RC2-1"

# ------------------------------------------------------------------
# 4) Clear out the old "4months" / "12months" / dividend-footnote cells —
#    that whole mini-block is replaced by the new Portfolio/SPY layout.
# ------------------------------------------------------------------
$ws.Range("I20").ClearContents()
$ws.Range("J20").ClearContents()
$ws.Range("F22").ClearContents()

# ------------------------------------------------------------------
# 5) New "Portfolio" column header next to the totals row.
# ------------------------------------------------------------------
$ws.Range("G19").Value = "Portfolio"
$ws.Range("G19").Font.Bold = $true
$ws.Range("G19").Font.Underline = $true

# ------------------------------------------------------------------
# 6) Footnote under the portfolio IRR block (replaces "12months").
# ------------------------------------------------------------------
$ws.Range("I21").Value = "(not including dividends from high yield stocks like OXY)"

# ------------------------------------------------------------------
# 7) SPY benchmark block (rows 23-25): header, mini Date/SPY price table,
#    and the %% / IRR formulas mirroring the portfolio's G20:H21 pattern.
# ------------------------------------------------------------------
$ws.Range("G23").Value = "SPY (benchmark)"
$ws.Range("G23").Font.Bold = $true
$ws.Range("G23").Font.Underline = $true

$ws.Range("J23").Value = "Date"
$ws.Range("K23").Value = "SPY"

$ws.Range("G24").Value = "%"
$ws.Range("H24").Formula = "=K25/K24-1"
$ws.Range("J24").Value = 45757
$ws.Range("K24").Value = 524.5
$ws.Range("K24").NumberFormat = "#,##0.00"

$ws.Range("G25").Value = "shares"
$ws.Range("H25").Formula = "=(H24+1)^3-1"
$ws.Range("J25").Value = 45868
$ws.Range("K25").Value = 634.4
$ws.Range("K25").NumberFormat = "#,##0.00"

# ------------------------------------------------------------------
# 8) Retrospective Analysis callout (rows 27-29).
# ------------------------------------------------------------------
$ws.Range("G27").Value = "Retrospective Analysis"
$ws.Range("G27").Font.Bold = $true
$ws.Range("G27").Font.Underline = $true

$ws.Range("G28").Value = "mock portfolio performed ~30% better than benchmark (SPY)"
$ws.Range("G29").Value = "key winners were PYPL and SMCI both were deep value investments trading at low ratios"

# ------------------------------------------------------------------
# 9) Window/view tidy-up: scroll the frozen pane down and land the
#    selection near the bottom of the new content.
# ------------------------------------------------------------------
$ws.Range("B8").Select()
$ws.Range("F22").Select()
